$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 10.2210306930377
$ws.Cells.Item(2, 3).Value = 11.51775305135823
$ws.Cells.Item(2, 4).Value = 5.235352649440111
$ws.Cells.Item(2, 6).Value = 24.23461331518008
$ws.Cells.Item(2, 7).Value = 28.6255903374972
$ws.Cells.Item(2, 8).Value = 14.38510219455232
$ws.Cells.Item(2, 9).Value = 23.26486237457069
$ws.Cells.Item(2, 11).Value = 8.58116882421815
$ws.Cells.Item(2, 12).Value = 11.21595893292223
$ws.Cells.Item(2, 13).Value = 13.75308397553642
$ws.Cells.Item(2, 15).Value = 21.85426183704564

# Row 3
$ws.Cells.Item(3, 2).Value = 9.923992341463785
$ws.Cells.Item(3, 3).Value = 11.53009331592607
$ws.Cells.Item(3, 4).Value = 5.154982484669983
$ws.Cells.Item(3, 6).Value = 24.29612488448738
$ws.Cells.Item(3, 7).Value = 28.73442185032179
$ws.Cells.Item(3, 8).Value = 14.43290714826732
$ws.Cells.Item(3, 9).Value = 23.36380770168875
$ws.Cells.Item(3, 11).Value = 8.337340936199
$ws.Cells.Item(3, 12).Value = 11.22606745934483
$ws.Cells.Item(3, 13).Value = 13.70532880120398
$ws.Cells.Item(3, 15).Value = 21.93734974384658

# Row 4
$ws.Cells.Item(4, 2).Value = 9.738061449572168
$ws.Cells.Item(4, 3).Value = 11.53844143672482
$ws.Cells.Item(4, 4).Value = 5.10429642810818
$ws.Cells.Item(4, 6).Value = 24.33981841187791
$ws.Cells.Item(4, 7).Value = 28.80958767732322
$ws.Cells.Item(4, 8).Value = 14.46429970509273
$ws.Cells.Item(4, 9).Value = 23.42859982846176
$ws.Cells.Item(4, 11).Value = 8.182608094998681
$ws.Cells.Item(4, 12).Value = 11.2340243499668
$ws.Cells.Item(4, 13).Value = 13.6777560088864
$ws.Cells.Item(4, 15).Value = 21.99255363269664

# Row 5
$ws.Cells.Item(5, 2).Value = 9.661512211389804
$ws.Cells.Item(5, 3).Value = 11.54203786455568
$ws.Cells.Item(5, 4).Value = 5.083317530828728
$ws.Cells.Item(5, 6).Value = 24.35911089421014
$ws.Cells.Item(5, 7).Value = 28.84230870328967
$ws.Cells.Item(5, 8).Value = 14.47760585606532
$ws.Cells.Item(5, 9).Value = 23.45601927652849
$ws.Cells.Item(5, 11).Value = 8.118350140601684
$ws.Cells.Item(5, 12).Value = 11.23770743687998
$ws.Cells.Item(5, 13).Value = 13.66696743736208
$ws.Cells.Item(5, 15).Value = 22.01610188223956

# Row 6
$ws.Cells.Item(6, 2).Value = 9.648757507038264
$ws.Cells.Item(6, 3).Value = 11.54264681329564
$ws.Cells.Item(6, 4).Value = 5.07981481544165
$ws.Cells.Item(6, 6).Value = 24.36240411646467
$ws.Cells.Item(6, 7).Value = 28.84786802930375
$ws.Cells.Item(6, 8).Value = 14.47984635477465
$ws.Cells.Item(6, 9).Value = 23.46063363959223
$ws.Cells.Item(6, 11).Value = 8.107609334909313
$ws.Cells.Item(6, 12).Value = 11.23834563331829
$ws.Cells.Item(6, 13).Value = 13.66520325324483
$ws.Cells.Item(6, 15).Value = 22.02007557945049

# Row 7
$ws.Cells.Item(7, 2).Value = 9.737032088882197
$ws.Cells.Item(7, 3).Value = 11.53848915109958
$ws.Cells.Item(7, 4).Value = 5.104014793864618
$ws.Cells.Item(7, 6).Value = 24.34007258027344
$ws.Cells.Item(7, 7).Value = 28.81002051101904
$ws.Cells.Item(7, 8).Value = 14.46447707707435
$ws.Cells.Item(7, 9).Value = 23.42896550187798
$ws.Cells.Item(7, 11).Value = 8.181746278462244
$ws.Cells.Item(7, 12).Value = 11.23407223686735
$ws.Cells.Item(7, 13).Value = 13.67760868834023
$ws.Cells.Item(7, 15).Value = 21.99286695327871

# Row 8
$ws.Cells.Item(8, 2).Value = 10.11941269024777
$ws.Cells.Item(8, 3).Value = 11.52184827625713
$ws.Cells.Item(8, 4).Value = 5.207926253222304
$ws.Cells.Item(8, 6).Value = 24.25459107708086
$ws.Cells.Item(8, 7).Value = 28.6613794614524
$ws.Cells.Item(8, 8).Value = 14.40116217029727
$ws.Cells.Item(8, 9).Value = 23.29814061601718
$ws.Cells.Item(8, 11).Value = 8.498173422581678
$ws.Cells.Item(8, 12).Value = 11.21908141411737
$ws.Cells.Item(8, 13).Value = 13.73625973762245
$ws.Cells.Item(8, 15).Value = 21.88204094908405

# Row 9
$ws.Cells.Item(9, 2).Value = 10.83661249468108
$ws.Cells.Item(9, 3).Value = 11.49530849476104
$ws.Cells.Item(9, 4).Value = 5.400546257688555
$ws.Cells.Item(9, 6).Value = 24.13409307638958
$ws.Cells.Item(9, 7).Value = 28.43641327454115
$ws.Cells.Item(9, 8).Value = 14.29317011156167
$ws.Cells.Item(9, 9).Value = 23.07362436414085
$ws.Cells.Item(9, 11).Value = 9.076475660650519
$ws.Cells.Item(9, 12).Value = 11.20354727657848
$ws.Cells.Item(9, 13).Value = 13.86479542039905
$ws.Cells.Item(9, 15).Value = 21.69797256922764

# Row 10
$ws.Cells.Item(10, 2).Value = 11.33818404053756
$ws.Cells.Item(10, 3).Value = 11.47948958065632
$ws.Cells.Item(10, 4).Value = 5.534594648093402
$ws.Cells.Item(10, 6).Value = 24.074436041014
$ws.Cells.Item(10, 7).Value = 28.31209881036975
$ws.Cells.Item(10, 8).Value = 14.22365838330871
$ws.Cells.Item(10, 9).Value = 22.92817254122809
$ws.Cells.Item(10, 11).Value = 9.472795806809534
$ws.Cells.Item(10, 12).Value = 11.20054812274644
$ws.Cells.Item(10, 13).Value = 13.96697196097528
$ws.Cells.Item(10, 15).Value = 21.58306078886712

# Row 11
$ws.Cells.Item(11, 2).Value = 11.55986417388516
$ws.Cells.Item(11, 3).Value = 11.4730848574007
$ws.Cells.Item(11, 4).Value = 5.593823064634205
$ws.Cells.Item(11, 6).Value = 24.05358822976941
$ws.Cells.Item(11, 7).Value = 28.26452136780071
$ws.Cells.Item(11, 8).Value = 14.19416473897815
$ws.Cells.Item(11, 9).Value = 22.86623204693775
$ws.Cells.Item(11, 11).Value = 9.646385833941808
$ws.Cells.Item(11, 12).Value = 11.20100053861492
$ws.Cells.Item(11, 13).Value = 14.01502143575919
$ws.Cells.Item(11, 15).Value = 21.53520732738713

# Row 12
$ws.Cells.Item(12, 2).Value = 11.64280359815114
$ws.Cells.Item(12, 3).Value = 11.47077274957391
$ws.Cells.Item(12, 4).Value = 5.615989205667307
$ws.Cells.Item(12, 6).Value = 24.04659938269301
$ws.Cells.Item(12, 7).Value = 28.24780133731905
$ws.Cells.Item(12, 8).Value = 14.18330184028143
$ws.Cells.Item(12, 9).Value = 22.84338431712628
$ws.Cells.Item(12, 11).Value = 9.711120190162717
$ws.Cells.Item(12, 12).Value = 11.20143202575092
$ws.Cells.Item(12, 13).Value = 14.03343192691747
$ws.Cells.Item(12, 15).Value = 21.51772299953636

# Row 13
$ws.Cells.Item(13, 2).Value = 11.62498693905704
$ws.Cells.Item(13, 3).Value = 11.47126567700255
$ws.Cells.Item(13, 4).Value = 5.611227162487279
$ws.Cells.Item(13, 6).Value = 24.04806425728951
$ws.Cells.Item(13, 7).Value = 28.25134454079754
$ws.Cells.Item(13, 8).Value = 14.18562776732065
$ws.Cells.Item(13, 9).Value = 22.84827795155323
$ws.Cells.Item(13, 11).Value = 9.697223533084303
$ws.Cells.Item(13, 12).Value = 11.20132754473788
$ws.Cells.Item(13, 13).Value = 14.02945749592456
$ws.Cells.Item(13, 15).Value = 21.52146022620104

# Row 14
$ws.Cells.Item(14, 2).Value = 11.56670820536519
$ws.Cells.Item(14, 3).Value = 11.47289237307604
$ws.Cells.Item(14, 4).Value = 5.59565200346739
$ws.Cells.Item(14, 6).Value = 24.05299509487838
$ws.Cells.Item(14, 7).Value = 28.26311978371696
$ws.Cells.Item(14, 8).Value = 14.19326491633667
$ws.Cells.Item(14, 9).Value = 22.86434016895238
$ws.Cells.Item(14, 11).Value = 9.651731817495042
$ws.Cells.Item(14, 12).Value = 11.20103082979321
$ws.Cells.Item(14, 13).Value = 14.01653182566115
$ws.Cells.Item(14, 15).Value = 21.53375611387876

# Row 15
$ws.Cells.Item(15, 2).Value = 11.53087773843904
$ws.Cells.Item(15, 3).Value = 11.47390349905674
$ws.Cells.Item(15, 4).Value = 5.586077295031341
$ws.Cells.Item(15, 6).Value = 24.05613336101676
$ws.Cells.Item(15, 7).Value = 28.27050147633475
$ws.Cells.Item(15, 8).Value = 14.1979826953869
$ws.Cells.Item(15, 9).Value = 22.87425789559382
$ws.Cells.Item(15, 11).Value = 9.623735506598306
$ws.Cells.Item(15, 12).Value = 11.2008829306806
$ws.Cells.Item(15, 13).Value = 14.00864219088664
$ws.Cells.Item(15, 15).Value = 21.54137065415562

# Row 16
$ws.Cells.Item(16, 2).Value = 11.32355975665069
$ws.Cells.Item(16, 3).Value = 11.4799240196261
$ws.Cells.Item(16, 4).Value = 5.53068778044221
$ws.Cells.Item(16, 6).Value = 24.07592518790885
$ws.Cells.Item(16, 7).Value = 28.31538921673941
$ws.Cells.Item(16, 8).Value = 14.22562865408523
$ws.Cells.Item(16, 9).Value = 22.93230554574435
$ws.Cells.Item(16, 11).Value = 9.461313232595458
$ws.Cells.Item(16, 12).Value = 11.20055502441168
$ws.Cells.Item(16, 13).Value = 13.96386250171672
$ws.Cells.Item(16, 15).Value = 21.58627716916013

# Row 17
$ws.Cells.Item(17, 2).Value = 11.1946608989443
$ws.Cells.Item(17, 3).Value = 11.48381969481651
$ws.Cells.Item(17, 4).Value = 5.496251942747989
$ws.Cells.Item(17, 6).Value = 24.08967896991726
$ws.Cells.Item(17, 7).Value = 28.34522929893724
$ws.Cells.Item(17, 8).Value = 14.243133306735
$ws.Cells.Item(17, 9).Value = 22.96899837831219
$ws.Cells.Item(17, 11).Value = 9.359929113242181
$ws.Cells.Item(17, 12).Value = 11.20081844741188
$ws.Cells.Item(17, 13).Value = 13.93678612785043
$ws.Cells.Item(17, 15).Value = 21.61495882622754

# Row 18
$ws.Cells.Item(18, 2).Value = 11.11991495856569
$ws.Cells.Item(18, 3).Value = 11.486134907207
$ws.Cells.Item(18, 4).Value = 5.476281157438914
$ws.Cells.Item(18, 6).Value = 24.09818182619989
$ws.Cells.Item(18, 7).Value = 28.36323685192567
$ws.Cells.Item(18, 8).Value = 14.25340179248344
$ws.Cells.Item(18, 9).Value = 22.99050091581804
$ws.Cells.Item(18, 11).Value = 9.300987163600329
$ws.Cells.Item(18, 12).Value = 11.201141034246
$ws.Cells.Item(18, 13).Value = 13.92136082373829
$ws.Cells.Item(18, 15).Value = 21.63187176997391

# Row 19
$ws.Cells.Item(19, 2).Value = 11.0945054137451
$ws.Cells.Item(19, 3).Value = 11.48693161532673
$ws.Cells.Item(19, 4).Value = 5.469491523216258
$ws.Cells.Item(19, 6).Value = 24.10116238759232
$ws.Cells.Item(19, 7).Value = 28.36947872004
$ws.Cells.Item(19, 8).Value = 14.25691292981629
$ws.Cells.Item(19, 9).Value = 22.99784962529796
$ws.Cells.Item(19, 11).Value = 9.280923690742204
$ws.Cells.Item(19, 12).Value = 11.20127966943629
$ws.Cells.Item(19, 13).Value = 13.91616386527318
$ws.Cells.Item(19, 15).Value = 21.63766961966649

# Row 20
$ws.Cells.Item(20, 2).Value = 11.20844574998137
$ws.Cells.Item(20, 3).Value = 11.48339728486677
$ws.Cells.Item(20, 4).Value = 5.499934769384685
$ws.Cells.Item(20, 6).Value = 24.08815357182285
$ws.Cells.Item(20, 7).Value = 28.34196534083683
$ws.Cells.Item(20, 8).Value = 14.24124918089419
$ws.Cells.Item(20, 9).Value = 22.96505119348693
$ws.Cells.Item(20, 11).Value = 9.370786954356078
$ws.Cells.Item(20, 12).Value = 11.20077270761696
$ws.Cells.Item(20, 13).Value = 13.93965317896141
$ws.Cells.Item(20, 15).Value = 21.61186255119402

# Row 21
$ws.Cells.Item(21, 2).Value = 11.58385393192799
$ws.Cells.Item(21, 3).Value = 11.47241150480934
$ws.Cells.Item(21, 4).Value = 5.600234007643584
$ws.Cells.Item(21, 6).Value = 24.05152219870684
$ws.Cells.Item(21, 7).Value = 28.25962587645416
$ws.Cells.Item(21, 8).Value = 14.19101340568339
$ws.Cells.Item(21, 9).Value = 22.85960581042769
$ws.Cells.Item(21, 11).Value = 9.665121265315772
$ws.Cells.Item(21, 12).Value = 11.20111093050063
$ws.Cells.Item(21, 13).Value = 14.02032265225178
$ws.Cells.Item(21, 15).Value = 21.53012722416253

# Row 22
$ws.Cells.Item(22, 2).Value = 11.82331216617186
$ws.Cells.Item(22, 3).Value = 11.46589133725362
$ws.Cells.Item(22, 4).Value = 5.664251524034357
$ws.Cells.Item(22, 6).Value = 24.03286122322682
$ws.Cells.Item(22, 7).Value = 28.21337180815772
$ws.Cells.Item(22, 8).Value = 14.15996319669849
$ws.Cells.Item(22, 9).Value = 22.79423388132966
$ws.Cells.Item(22, 11).Value = 9.851640021286489
$ws.Cells.Item(22, 12).Value = 11.20284784580225
$ws.Cells.Item(22, 13).Value = 14.07429413313313
$ws.Cells.Item(22, 15).Value = 21.48042035699536

# Row 23
$ws.Cells.Item(23, 2).Value = 11.6960702600381
$ws.Cells.Item(23, 3).Value = 11.46931110426889
$ws.Cells.Item(23, 4).Value = 5.630227870736126
$ws.Cells.Item(23, 6).Value = 24.04233755717783
$ws.Cells.Item(23, 7).Value = 28.23736490094112
$ws.Cells.Item(23, 8).Value = 14.17637232202184
$ws.Cells.Item(23, 9).Value = 22.82879992197451
$ws.Cells.Item(23, 11).Value = 9.752637624346509
$ws.Cells.Item(23, 12).Value = 11.20178250704121
$ws.Cells.Item(23, 13).Value = 14.04537773452857
$ws.Cells.Item(23, 15).Value = 21.50660985653262

# Row 24
$ws.Cells.Item(24, 2).Value = 11.20221560749456
$ws.Cells.Item(24, 3).Value = 11.48358802107314
$ws.Cells.Item(24, 4).Value = 5.498270302020098
$ws.Cells.Item(24, 6).Value = 24.08884134908751
$ws.Cells.Item(24, 7).Value = 28.34343832256869
$ws.Cells.Item(24, 8).Value = 14.24210035614334
$ws.Cells.Item(24, 9).Value = 22.96683444691794
$ws.Cells.Item(24, 11).Value = 9.36588016131023
$ws.Cells.Item(24, 12).Value = 11.20079285341484
$ws.Cells.Item(24, 13).Value = 13.93835654434849
$ws.Cells.Item(24, 15).Value = 21.61326105811865

# Row 25
$ws.Cells.Item(25, 2).Value = 10.64668045431965
$ws.Cells.Item(25, 3).Value = 11.50183948600876
$ws.Cells.Item(25, 4).Value = 5.349702564857011
$ws.Cells.Item(25, 6).Value = 24.16162924847785
$ws.Cells.Item(25, 7).Value = 28.49010664840108
$ws.Cells.Item(25, 8).Value = 14.32065668677324
$ws.Cells.Item(25, 9).Value = 23.13093571724795
$ws.Cells.Item(25, 11).Value = 9.076475660650519
$ws.Cells.Item(25, 12).Value = 11.20354727657848
$ws.Cells.Item(25, 13).Value = 13.86479542039905
$ws.Cells.Item(25, 15).Value = 21.69797256922764

